$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers (human readable labels, re-keyed to new text/order)
$ws.Range("A1").Value = "Comarca nombre"
$ws.Range("B1").Value = "Número hogares"
$ws.Range("C1").Value = "Comarca código"
$ws.Range("D1").Value = "Provincia código"
$ws.Range("E1").Value = "Aragón"
$ws.Range("F1").Value = "Municipio código"
$ws.Range("G1").Value = "Provincia nombre"
$ws.Range("H1").Value = "Estructura hogar"
$ws.Range("I1").Value = "Municipio nombre"

# Row 2 - concept/measure URI row
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "iaest-measure:numero-hogares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "iaest-measure:estructura-hogar"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3 - role row (dim / medida / null)
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "dim"

# Row 4 - codelist/type row
$ws.Range("A4").Value = "URI-comarca"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-Provincia"
$ws.Range("H4").Value = "xsd:string"
$ws.Range("I4").Value = "URI-Municipio"
